# Inserts two new weekly price rows ("Fruta / hortaliza, semanal") into the
# Tomate subset sheet, just before the existing row for 2021-08-10
# (old row 396). This pushes all subsequent data rows down by two rows
# (old 396..489 -> new 398..491) and fills the two freshly inserted rows
# with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 396 downward by inserting two blank rows above the old row 396.
$ws.Rows("396:397").Insert()

# --- New row 396 ---
$ws.Cells.Item(396, 1).Value = 7
$ws.Cells.Item(396, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(396, 3).Value = "Ñuble"
$ws.Cells.Item(396, 4).Value = 44754
$ws.Cells.Item(396, 5).Value = 16
$ws.Cells.Item(396, 6).Value = 100112020
$ws.Cells.Item(396, 7).Value = "Tomate"
$ws.Cells.Item(396, 8).Value = "Larga vida"
$ws.Cells.Item(396, 9).Value = "Primera"
$ws.Cells.Item(396, 10).Value = 500
$ws.Cells.Item(396, 11).Value = 8000
$ws.Cells.Item(396, 12).Value = 9000
$ws.Cells.Item(396, 13).Value = 8500
$ws.Cells.Item(396, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(396, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(396, 16).Value = 472
$ws.Cells.Item(396, 17).Value = 18
$ws.Cells.Item(396, 18).Value = "Hortaliza"

# --- New row 397 ---
$ws.Cells.Item(397, 1).Value = 7
$ws.Cells.Item(397, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(397, 3).Value = "Ñuble"
$ws.Cells.Item(397, 4).Value = 44754
$ws.Cells.Item(397, 5).Value = 16
$ws.Cells.Item(397, 6).Value = 100112020
$ws.Cells.Item(397, 7).Value = "Tomate"
$ws.Cells.Item(397, 8).Value = "Larga vida"
$ws.Cells.Item(397, 9).Value = "Primera"
$ws.Cells.Item(397, 10).Value = 600
$ws.Cells.Item(397, 11).Value = 3500
$ws.Cells.Item(397, 12).Value = 4000
$ws.Cells.Item(397, 13).Value = 3750
$ws.Cells.Item(397, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(397, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(397, 16).Value = 375
$ws.Cells.Item(397, 17).Value = 10
$ws.Cells.Item(397, 18).Value = "Hortaliza"
